$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the label in A10 from "L_uF" to "L_uH"
$ws.Range("A10").Value = "L_uH"

# Move the active selection to A11, matching the author's final cursor position
$ws.Range("A11").Select()
